$d = $word.ActiveDocument

# 1) "Tamanho" -> "Descrição" (single occurrence, in the "Tipos de Marmitas" attribute list)
$d.Content.Find.Execute("Tamanho", $true, $false, $false, $false, $false, $true, 1, $false, "Descrição", 2) | Out-Null

# 2) Add a new "Valor por Peso" bullet right after the "Observação" bullet in the
#    "Marmitas" attributes list (same list style / numbering as its siblings).
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Observação*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $obsPara = $d.Paragraphs.Item($targetIndex)
    $obsPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "Valor por Peso"
    $newPara.Style = $obsPara.Style
    $newPara.Range.ListFormat.ApplyListTemplateWithLevel($obsPara.Range.ListFormat.ListTemplate, $true)
    $newPara.Alignment = $obsPara.Alignment
}
